$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a "last changed" date that is bumped by
# one day (serial 46081 -> 46082, i.e. 2026-02-28 -> 2026-03-01) for
# every data row (rows 2 through 537).
$lastRow = 537
$range = $ws.Range("C2:C$lastRow")
$range.Value = 46082
